{"js": "// Update the date paragraph (first paragraph in the body) and the\n// division-problem table cells. Each cell's old text is unique-ish in the\n// document but a few NEW values collide with OTHER cells' OLD values\n// (e.g. \"79\u00f75=\" and \"18\u00f78=\" each occur twice across old/new), so we must\n// address every paragraph/cell by its *position* in the document rather\n// than doing a blind global find-and-replace (which could double-replace\n// or replace the wrong occurrence depending on execution order).\n\nconst body = context.document.body;\n\n// --- 1) Date paragraph -----------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-05-19 Sunday\", \"Replace\");\n\n// --- 2) Table of division problems ------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Row-major list of [oldText, newText] for the 5x5 grid of populated\n// cells (the remaining rows in the table are blank spacer rows).\nconst cellUpdates = [\n  [\"83\u00f74=\", \"39\u00f73=\"],\n  [\"56\u00f77=\", \"18\u00f78=\"],\n  [\"33\u00f75=\", \"41\u00f74=\"],\n  [\"79\u00f78=\", \"29\u00f75=\"],\n  [\"94\u00f77=\", \"38\u00f78=\"],\n\n  [\"32\u00f73=\", \"38\u00f74=\"],\n  [\"18\u00f77=\", \"13\u00f73=\"],\n  [\"54\u00f75=\", \"79\u00f75=\"],\n  [\"18\u00f72=\", \"34\u00f78=\"],\n  [\"22\u00f76=\", \"31\u00f73=\"],\n\n  [\"40\u00f78=\", \"30\u00f79=\"],\n  [\"70\u00f75=\", \"79\u00f73=\"],\n  [\"20\u00f73=\", \"27\u00f79=\"],\n  [\"23\u00f76=\", \"62\u00f73=\"],\n  [\"29\u00f76=\", \"18\u00f79=\"],\n\n  [\"60\u00f79=\", \"64\u00f79=\"],\n  [\"84\u00f77=\", \"22\u00f72=\"],\n  [\"11\u00f72=\", \"64\u00f76=\"],\n  [\"66\u00f75=\", \"57\u00f79=\"],\n  [\"33\u00f76=\", \"18\u00f78=\"],\n\n  [\"90\u00f77=\", \"98\u00f77=\"],\n  [\"74\u00f73=\", \"73\u00f76=\"],\n  [\"79\u00f75=\", \"37\u00f74=\"],\n  [\"40\u00f72=\", \"84\u00f78=\"],\n  [\"52\u00f74=\", \"44\u00f73=\"],\n];\n\nconst cols = 5;\n// Content lives in table rows 0, 4, 8, 12, 16 (every 4th row; the rows in\n// between are empty spacer rows), 5 cells per row.\nconst contentRows = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < cellUpdates.length; i++) {\n  const rowIdx = contentRows[Math.floor(i / cols)];\n  const colIdx = i % cols;\n  const cell = table.getCell(rowIdx, colIdx);\n  cell.value = cellUpdates[i][1];\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the division-problem table cells.\n#\n# Each cell's old text is updated to a new value, but several NEW values\n# collide with OTHER cells' OLD values (e.g. \"79\u00f75=\" and \"18\u00f78=\" each show\n# up twice across the old/new sets). A blind document-wide Find/Replace\n# pass (Content.Find.Execute with Replace:=wdReplaceAll) risks re-matching\n# an already-updated cell on a later replacement, or matching the wrong\n# occurrence first. So address every paragraph/cell by its position in the\n# document (first paragraph, then table row/column) instead.\n\n$d = $word.ActiveDocument\n\n# --- 1) Date paragraph --------------------------------------------------\n$d.Paragraphs.Item(1).Range.Text = \"2024-05-19 Sunday\"\n\n# --- 2) Table of division problems --------------------------------------\n$tbl = $d.Tables.Item(1)\n\n# Row-major list of new cell text for the 5x5 grid of populated cells (the\n# other rows in the 20-row table are blank spacer rows and are untouched).\n$newValues = @(\n    \"39\u00f73=\", \"18\u00f78=\", \"41\u00f74=\", \"29\u00f75=\", \"38\u00f78=\",\n    \"38\u00f74=\", \"13\u00f73=\", \"79\u00f75=\", \"34\u00f78=\", \"31\u00f73=\",\n    \"30\u00f79=\", \"79\u00f73=\", \"27\u00f79=\", \"62\u00f73=\", \"18\u00f79=\",\n    \"64\u00f79=\", \"22\u00f72=\", \"64\u00f76=\", \"57\u00f79=\", \"18\u00f78=\",\n    \"98\u00f77=\", \"73\u00f76=\", \"37\u00f74=\", \"84\u00f78=\", \"44\u00f73=\"\n)\n\n# Content lives in table rows 1, 5, 9, 13, 17 (1-based; every 4th row, the\n# rows in between are empty spacer rows), 5 cells per row.\n$contentRows = @(1, 5, 9, 13, 17)\n$cols = 5\n\nfor ($i = 0; $i -lt $newValues.Length; $i++) {\n    $rowIdx = $contentRows[[math]::Floor($i / $cols)]\n    $colIdx = ($i % $cols) + 1\n    $tbl.Cell($rowIdx, $colIdx).Range.Text = $newValues[$i]\n}\n"}
